$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.400.64"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "3.459.54"
$ws.Range("E3").Value = "  -1.75%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'588.04"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "'176.64"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").Value = "'0.618"
$ws.Range("E7").Value = "  +3.03%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "3.460.38"
$ws.Range("E9").Value = "  -1.76%  "
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").Value = "'0.418"
$ws.Range("E12").Value = "  -1.67%  "
$ws.Range("D13").Value = "4.060.34"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("D15").Value = "'30.08"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").Value = "66.312.51"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").Value = "3.451.46"
$ws.Range("E18").Value = "  -2.09%  "
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("D20").Value = "'13.84"
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").Value = "'372.63"
$ws.Range("E21").Value = "  -2.97%  "
$ws.Range("D22").Value = "'7.62"
$ws.Range("E22").Value = "  -3.01%  "
$ws.Range("D23").Value = "'73.20"
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("E25").Value = "  +3.63%  "
$ws.Range("D27").Value = "'9.94"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  +2.42%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  -0.90%  "
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("E32").Value = "  -3.29%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E34").Value = "  -3.01%  "
$ws.Range("E35").Value = "  -7.04%  "
$ws.Range("E36").Value = "  -2.64%  "
$ws.Range("D37").Value = "'160.84"
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("D39").Value = "'28.22"
$ws.Range("E39").Value = "  -4.62%  "
$ws.Range("D40").Value = "'1.81"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").Value = "2.769.56"
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("D44").Value = "'6.47"
$ws.Range("E44").Value = "  -2.33%  "
$ws.Range("D45").Value = "'0.0694"
$ws.Range("E45").Value = "  -1.84%  "
$ws.Range("D46").Value = "'25.34"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").Value = "'337.97"
$ws.Range("E47").Value = "  +3.66%  "
$ws.Range("D48").Value = "'40.01"
$ws.Range("E48").Value = "  -1.75%  "
$ws.Range("E49").Value = "  -1.75%  "
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -2.64%  "
